$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell F1 "time_taken", copying the header formatting from E1
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("F1").Value = "time_taken"

# Timestamps for F2:F76 (time_taken values per row)
$timeTaken = @(
    "2021-10-05 10:52:14.787024",
    "2021-10-05 10:52:14.787035",
    "2021-10-05 10:52:14.787039",
    "2021-10-05 10:52:14.787041",
    "2021-10-05 10:52:14.787044",
    "2021-10-05 10:52:14.787047",
    "2021-10-05 10:52:14.787050",
    "2021-10-05 10:52:14.787052",
    "2021-10-05 10:52:14.787055",
    "2021-10-05 10:52:14.787057",
    "2021-10-05 10:52:14.787060",
    "2021-10-05 10:52:14.787062",
    "2021-10-05 10:52:14.787065",
    "2021-10-05 10:52:14.787068",
    "2021-10-05 10:52:14.787070",
    "2021-10-05 10:52:14.787072",
    "2021-10-05 10:52:14.787075",
    "2021-10-05 10:52:14.787078",
    "2021-10-05 10:52:14.787081",
    "2021-10-05 10:52:14.787083",
    "2021-10-05 10:52:14.787086",
    "2021-10-05 10:52:14.787088",
    "2021-10-05 10:52:14.787091",
    "2021-10-05 10:52:14.787093",
    "2021-10-05 10:52:14.787096",
    "2021-10-05 10:52:14.787099",
    "2021-10-05 10:52:14.787101",
    "2021-10-05 10:52:14.787104",
    "2021-10-05 10:52:14.787106",
    "2021-10-05 10:52:14.787109",
    "2021-10-05 10:52:14.787111",
    "2021-10-05 10:52:14.787113",
    "2021-10-05 10:52:14.787116",
    "2021-10-05 10:52:14.787119",
    "2021-10-05 10:52:14.787122",
    "2021-10-05 10:52:14.787124",
    "2021-10-05 10:52:14.787126",
    "2021-10-05 10:52:14.787129",
    "2021-10-05 10:52:14.787131",
    "2021-10-05 10:52:14.787134",
    "2021-10-05 10:52:14.787137",
    "2021-10-05 10:52:14.787140",
    "2021-10-05 10:52:14.787142",
    "2021-10-05 10:52:14.787145",
    "2021-10-05 10:52:14.787147",
    "2021-10-05 10:52:14.787150",
    "2021-10-05 10:52:14.787152",
    "2021-10-05 10:52:14.787154",
    "2021-10-05 10:52:14.787157",
    "2021-10-05 10:52:14.787159",
    "2021-10-05 10:52:14.787162",
    "2021-10-05 10:52:14.787164",
    "2021-10-05 10:52:14.787167",
    "2021-10-05 10:52:14.787170",
    "2021-10-05 10:52:14.787172",
    "2021-10-05 10:52:14.787174",
    "2021-10-05 10:52:14.787177",
    "2021-10-05 10:52:14.787179",
    "2021-10-05 10:52:14.787182",
    "2021-10-05 10:52:14.787184",
    "2021-10-05 10:52:14.787187",
    "2021-10-05 10:52:14.787189",
    "2021-10-05 10:52:14.787192",
    "2021-10-05 10:52:14.787194",
    "2021-10-05 10:52:14.787198",
    "2021-10-05 10:52:14.787201",
    "2021-10-05 10:52:14.787203",
    "2021-10-05 10:52:14.787206",
    "2021-10-05 10:52:14.787208",
    "2021-10-05 10:52:14.787211",
    "2021-10-05 10:52:14.787213",
    "2021-10-05 10:52:14.787216",
    "2021-10-05 10:52:14.787218",
    "2021-10-05 10:52:14.787221",
    "2021-10-05 10:52:14.787224"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timeTaken[$i]
}
